$wb = $excel.ActiveWorkbook

# Sheet "建物" (building): property_category column (I) for all data rows (2-9)
# was incorrectly set to "land" ("land"); fix it to "building" ("building")
$wsBuilding = $wb.Worksheets.Item("建物")
for ($row = 2; $row -le 9; $row++) {
    $wsBuilding.Cells.Item($row, 9).Value = "building"
}

# Sheet "汽車" (car): property_category column (H) for the data row (2)
# was incorrectly set to "land" ("land"); fix it to "car" ("car")
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value = "car"
